# Auto-generated edit script: apply cyclic permutation of species-observation blocks
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{}
$rows[3] = [ordered]@{
    A = 111866265
    B = 78107
    E = 6453
    F = "Vedskivlav"
    G = "Hertelidea botryosa"
    H = "(Fr.) Printzen & Kantvilas"
    Q = 702680.6244306123
    R = 7299924.914052285
}
$rows[4] = [ordered]@{
    A = 111866048
    B = 90682
    E = 2059
    F = "Skrovlig taggsvamp"
    G = "Hydnellum scabrosum"
    H = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
    Q = 702750.1350314748
    R = 7299799.924799141
}
$rows[5] = [ordered]@{
    A = 111866159
    B = 90652
    E = 3100
    F = "Talltaggsvamp"
    G = "Bankera fuligineoalba"
    H = "(Schmidt : Fr.) Pouzar"
    Q = 702755.4455659754
    R = 7299865.042498757
}
$rows[6] = [ordered]@{
    A = 111866301
    B = 90660
    E = 4362
    F = "Blå taggsvamp"
    G = "Hydnellum caeruleum"
    H = "(Hornem.) P.Karst."
    Q = 702522.1051459431
    R = 7300047.742725079
}
$rows[7] = [ordered]@{
    A = 111865668
    B = 78107
    E = 6453
    F = "Vedskivlav"
    G = "Hertelidea botryosa"
    H = "(Fr.) Printzen & Kantvilas"
    Q = 702740.9003275807
    R = 7299743.601162716
}
$rows[8] = [ordered]@{
    A = 111866031
    B = 78107
    E = 6453
    F = "Vedskivlav"
    G = "Hertelidea botryosa"
    H = "(Fr.) Printzen & Kantvilas"
    Q = 702750.1350314748
    R = 7299799.924799141
}
$rows[9] = [ordered]@{
    A = 111865263
    B = 90658
    E = 4361
    F = "Orange taggsvamp"
    G = "Hydnellum aurantiacum"
    H = "(Batsch:Fr.) P.Karst."
    Q = 702714.1819675351
    R = 7299724.394724619
}
$rows[10] = [ordered]@{
    A = 111866021
    B = 78107
    E = 6453
    F = "Vedskivlav"
    G = "Hertelidea botryosa"
    H = "(Fr.) Printzen & Kantvilas"
    Q = 702738.1111920479
    R = 7299806.49869829
}
$rows[12] = [ordered]@{
    A = 111866194
    B = 90682
    E = 2059
    F = "Skrovlig taggsvamp"
    G = "Hydnellum scabrosum"
    H = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
    Q = 702686.7518818546
    R = 7299919.985876646
}
$rows[13] = [ordered]@{
    A = 111865866
    B = 90652
    E = 3100
    F = "Talltaggsvamp"
    G = "Bankera fuligineoalba"
    H = "(Schmidt : Fr.) Pouzar"
    Q = 702753.3055412351
    R = 7299801.798166115
}
$rows[14] = [ordered]@{
    A = 111866170
    B = 90682
    E = 2059
    F = "Skrovlig taggsvamp"
    G = "Hydnellum scabrosum"
    H = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
    Q = 702754.3208386695
    R = 7299886.818591502
}
$rows[16] = [ordered]@{
    A = 111865961
    B = 77267
    E = 6446
    F = "Kolflarnlav"
    G = "Carbonicola anthracophila"
    H = "(Nyl.) Bendiksby & Timdal"
    Q = 702714.4770808229
    R = 7299790.39698876
}
$rows[18] = [ordered]@{
    A = 111866065
    B = 78107
    E = 6453
    F = "Vedskivlav"
    G = "Hertelidea botryosa"
    H = "(Fr.) Printzen & Kantvilas"
    Q = 702767.9701038125
    R = 7299827.988589783
}
$rows[19] = [ordered]@{
    A = 111865578
    B = 90854
    E = 2079
    F = "Nordtagging"
    G = "Odonticium romellii"
    H = "(S.Lundell) Parmasto"
    Q = 702741.9879008483
    R = 7299745.739876431
}
$rows[20] = [ordered]@{
    A = 111866276
    B = 78107
    E = 6453
    F = "Vedskivlav"
    G = "Hertelidea botryosa"
    H = "(Fr.) Printzen & Kantvilas"
    Q = 702660.5304515015
    R = 7299928.856484808
}
$rows[21] = [ordered]@{
    A = 111865981
    B = 90652
    E = 3100
    F = "Talltaggsvamp"
    G = "Bankera fuligineoalba"
    H = "(Schmidt : Fr.) Pouzar"
    Q = 702695.6801449896
    R = 7299770.100652335
}
$rows[22] = [ordered]@{
    A = 111865524
    B = 90660
    E = 4362
    F = "Blå taggsvamp"
    G = "Hydnellum caeruleum"
    H = "(Hornem.) P.Karst."
    Q = 702731.0699128226
    R = 7299742.494774668
}

foreach ($r in $rows.Keys) {
    $v = $rows[$r]
    $ws.Range("A$r").Value = $v.A
    $ws.Range("B$r").Value = $v.B
    $ws.Range("E$r").Value = $v.E
    $ws.Range("F$r").Value = $v.F
    $ws.Range("G$r").Value = $v.G
    $ws.Range("H$r").Value = $v.H
    $ws.Range("Q$r").Value = $v.Q
    $ws.Range("R$r").Value = $v.R
}

# Move the public-comment text (AC) along with its observation block: row 5 -> row 12
$ws.Range("AC12").Value = "Flera fruktkoppar som växer i en häxring"
$ws.Range("AC5").Value = $null
